# Update "Elapsed Duration(Hrs)" (column G) values on several sheets to
# reflect the latest PCM refresh timestamp.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3926:13:07"
$ws.Range("G3").Value = "65:45:45"

$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12107:36:47"
$ws.Range("G3").Value = "3237:20:16"
$ws.Range("G4").Value = "475:31:50"

$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2953:26:36"
$ws.Range("G3").Value = "180:38:51"

$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "427:25:35"

$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "67:57:53"
